$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1616541353383459
$ws.Range("C2").Value = 0.6165413533834586
$ws.Range("J2").Value = 0.007518796992481203
$ws.Range("P2").Value = 0.1466165413533835
$ws.Range("S2").Value = 0.06766917293233082

# Row 3
$ws.Range("B3").Value = 0.0053475935828877
$ws.Range("C3").Value = 0.0106951871657754
$ws.Range("J3").Value = 0.0427807486631016
$ws.Range("P3").Value = 0.7593582887700535
$ws.Range("S3").Value = 0.1818181818181818

# Row 4
$ws.Range("J4").Value = 0.01818181818181818
$ws.Range("P4").Value = 0.7454545454545455
$ws.Range("S4").Value = 0.2363636363636364

# Row 6
$ws.Range("B6").Value = 0.06986899563318777
$ws.Range("D6").Value = 0.02620087336244541
$ws.Range("F6").Value = 0.04366812227074236
$ws.Range("J6").Value = 0.240174672489083
$ws.Range("O6").Value = 0.02183406113537118
$ws.Range("Q6").Value = 0.1834061135371179
$ws.Range("R6").Value = 0.07423580786026202
$ws.Range("S6").Value = 0.3406113537117904

# Row 7
$ws.Range("B7").Value = 0.09625668449197861
$ws.Range("D7").Value = 0.0481283422459893
$ws.Range("F7").Value = 0.0855614973262032
$ws.Range("J7").Value = 0.08021390374331551
$ws.Range("O7").Value = 0.0106951871657754
$ws.Range("Q7").Value = 0.2352941176470588
$ws.Range("S7").Value = 0.3529411764705883

# Row 8
$ws.Range("B8").Value = 0.1044776119402985
$ws.Range("D8").Value = 0.02132196162046908
$ws.Range("F8").Value = 0.07036247334754797
$ws.Range("J8").Value = 0.1023454157782516
$ws.Range("O8").Value = 0.01918976545842218
$ws.Range("Q8").Value = 0.2004264392324094
$ws.Range("R8").Value = 0.1023454157782516
$ws.Range("S8").Value = 0.3795309168443497

# Row 9
$ws.Range("B9").Value = 0.1488095238095238
$ws.Range("D9").Value = 0.0119047619047619
$ws.Range("E9").Value = 0.005952380952380952
$ws.Range("F9").Value = 0.07738095238095238
$ws.Range("J9").Value = 0.07142857142857142
$ws.Range("O9").Value = 0.0119047619047619
$ws.Range("Q9").Value = 0.2559523809523809
$ws.Range("R9").Value = 0.04166666666666666
$ws.Range("S9").Value = 0.375

# Row 10
$ws.Range("B10").Value = 0.0853950518754988
$ws.Range("D10").Value = 0.02154828411811652
$ws.Range("F10").Value = 0.07501995211492418
$ws.Range("J10").Value = 0.09976057462090981
$ws.Range("O10").Value = 0.0223463687150838
$ws.Range("Q10").Value = 0.2250598563447725
$ws.Range("R10").Value = 0.1021548284118117
$ws.Range("S10").Value = 0.3687150837988827

# Row 11
$ws.Range("G11").Value = 0.1293706293706294
$ws.Range("J11").Value = 0.1118881118881119
$ws.Range("K11").Value = 0.1783216783216783
$ws.Range("L11").Value = 0.5734265734265734
$ws.Range("S11").Value = 0.006993006993006993

# Row 12
$ws.Range("G12").Value = 0.7529411764705882
$ws.Range("J12").Value = 0.1764705882352941
$ws.Range("K12").Value = 0.005882352941176471
$ws.Range("L12").Value = 0.04117647058823529
$ws.Range("S12").Value = 0.02352941176470588

# Row 13
$ws.Range("G13").Value = 0.7741935483870968
$ws.Range("J13").Value = 0.1612903225806452
$ws.Range("S13").Value = 0.06451612903225806

# Row 14
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25

# Row 15
$ws.Range("F15").Value = 0.012
$ws.Range("H15").Value = 0.192
$ws.Range("I15").Value = 0.06
$ws.Range("J15").Value = 0.324
$ws.Range("K15").Value = 0.06
$ws.Range("M15").Value = 0.008
$ws.Range("O15").Value = 0.052
$ws.Range("S15").Value = 0.292

# Row 16
$ws.Range("F16").Value = 0.02336448598130841
$ws.Range("H16").Value = 0.1542056074766355
$ws.Range("I16").Value = 0.07943925233644859
$ws.Range("J16").Value = 0.3831775700934579
$ws.Range("K16").Value = 0.1355140186915888
$ws.Range("M16").Value = 0.02336448598130841
$ws.Range("O16").Value = 0.06074766355140187
$ws.Range("S16").Value = 0.1401869158878505

# Row 17
$ws.Range("F17").Value = 0.017578125
$ws.Range("H17").Value = 0.197265625
$ws.Range("I17").Value = 0.068359375
$ws.Range("J17").Value = 0.423828125
$ws.Range("K17").Value = 0.087890625
$ws.Range("M17").Value = 0.017578125
$ws.Range("N17").Value = 0.00390625
$ws.Range("O17").Value = 0.08984375
$ws.Range("S17").Value = 0.09375

# Row 18
$ws.Range("F18").Value = 0.02727272727272727
$ws.Range("H18").Value = 0.2090909090909091
$ws.Range("I18").Value = 0.05
$ws.Range("J18").Value = 0.4136363636363636
$ws.Range("K18").Value = 0.1045454545454545
$ws.Range("M18").Value = 0.00909090909090909
$ws.Range("O18").Value = 0.05
$ws.Range("S18").Value = 0.1363636363636364

# Row 19
$ws.Range("F19").Value = 0.01423785594639866
$ws.Range("H19").Value = 0.202680067001675
$ws.Range("I19").Value = 0.07705192629815745
$ws.Range("J19").Value = 0.4003350083752094
$ws.Range("K19").Value = 0.102177554438861
$ws.Range("M19").Value = 0.01172529313232831
$ws.Range("N19").Value = 0.002512562814070352
$ws.Range("O19").Value = 0.08040201005025126
$ws.Range("S19").Value = 0.1088777219430486
